$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (shifts existing rows 18-34 down to 19-35)
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with a new weekly record, matching the
# surrounding rows' fixed columns (Mercado/Region/Producto/etc.)
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 45090
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = "Tropicales y subtropicales"
$ws.Range("I18").Value = 100108001
$ws.Range("J18").Value = "Guayaba"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 2600
$ws.Range("O18").Value = 2600
$ws.Range("P18").Value = 2600
$ws.Range("Q18").Value = "$/kilo"
$ws.Range("R18").Value = "Región de Arica y Parinacota"
$ws.Range("S18").Value = 2600
$ws.Range("T18").Value = 1
